# Update "想去人数" (column F) counts on the 展览 and 全部类型 sheets
# to reflect the latest scrape output (gh-pages regeneration at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - row => new F value
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionUpdates = @{
    2  = 127
    3  = 317
    5  = 1236
    6  = 76
    7  = 2140
    11 = 4767
    14 = 296
    15 = 219
    16 = 24
    20 = 109
    21 = 3689
    22 = 463
    23 = 602
    24 = 23
    26 = 98
    27 = 110
    29 = 10
    32 = 6
    34 = 828
    35 = 2296
}
foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Sheet "全部类型" (all types) - row => new F value
$sheetAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2  = 127
    3  = 317
    5  = 1236
    6  = 76
    7  = 2140
    11 = 4767
    14 = 296
    15 = 219
    16 = 24
    20 = 109
    21 = 3689
    22 = 463
    23 = 602
    24 = 23
    26 = 98
    27 = 110
    29 = 10
    32 = 6
    35 = 828
    36 = 2297
}
foreach ($row in $allUpdates.Keys) {
    $sheetAll.Range("F$row").Value = $allUpdates[$row]
}
